$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the three "spell-checked word" paragraphs that had their runs
#    split around <w:proofErr> markers back into their natural prose, which
#    also removes the now pointless <w:proofErr> elements.
# ---------------------------------------------------------------------------

function Collapse-ParagraphText($fullText) {
    # Searching for the paragraph's own (unchanged) rendered text and
    # "replacing" it with the identical string forces the engine to
    # rebuild the matched range as a single run, which drops the now
    # pointless <w:proofErr> spell-check markers that used to split it.
    $found = $d.Content.Find.Execute($fullText, $true, $false, $false, $false,
                                      $false, $true, 1, $false, $fullText, 2)
    return $found
}

$vlcFull = "Gekozen om scherm op te nemen via VLC Media Player. Geluid wordt door een Lumix camera opgenomen tijdens het opnemen van de schermopnames en later in Premiere Pro samengezet met de schermopnames."
Collapse-ParagraphText $vlcFull | Out-Null

$vensterFull = "Venster modussen in GIMP en PS. Achter gekomen dat er ook veel liefde is voor de multi venster modus en dat PS 5 voor gedefinieerde modussen heeft voor elk van zijn soort gebruikers."
Collapse-ParagraphText $vensterFull | Out-Null

$pluginFull = "Plugin optie bij GIMP bespreken"
Collapse-ParagraphText $pluginFull | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a new list item after "Video flow bewerken." describing the video
#    batch edit, and move the hidden _GoBack bookmark so it still marks the
#    end of the document (now inside the newly added paragraph).
# ---------------------------------------------------------------------------

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newParaIndex)
$insertionStart = $newPara.Range.Start
$insertPoint = $d.Range($insertionStart, $insertionStart)
$insertPoint.InsertAfter("Video batch in filmpje plakken en knippen waar nodigX")

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$newParaIndex2 = $d.Paragraphs.Count
$newPara2 = $d.Paragraphs($newParaIndex2)
$endPos = $newPara2.Range.End - 1
$bookmarkRange = $d.Range($endPos - 1, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderRange = $d.Range($endPos - 1, $endPos)
$placeholderRange.Delete()
